$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix "unadited" -> "unaudited" typo in the two header cells.
$ws.Range("B1").Value = "Full Year (unaudited)"
$ws.Range("E1").Value = "2013 (unaudited)"

# Update window size for the workbook view.
$excel.Width = 9255
$excel.Height = 2595

# Update the active sheet view: scroll so column B is the left-most visible
# column, and move the selection to J4.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J4").Select()
